$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B2: shared string content change (kernel file -> notebook name) ---
$ws.Range("B2").Value = "TestJupyterConnection.ipynb"

# --- New G column: width 15, new cell G3 with xloJpyRun formula ---
$ws.Range("G1").ColumnWidth = 14.1666666
$ws.Range("G3").Formula = '=_xll.xloJpyRun(B3,"func2({})", 7)'

# --- B5:D7 watch array: swap argument order (B3 first, then label) ---
$ws.Range("B5:D7").FormulaArray = '=_xll.xloJpyWatch(B3,"watch_var")'

# --- B9: jptest formula kept (re-enter to refresh), now volatile in source ---
$ws.Range("B9").Formula = '=_xll.jptest(1)'

# --- remove row 10 (xloLog call no longer present) ---
$ws.Range("B10").ClearContents()

# --- new E13 cell: xloPyDebug call ---
$ws.Range("E13").Formula = '=_xll.xloPyDebug("pdb")'

# --- selection moves to F11 ---
$ws.Range("F11").Select()
